$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Intro paragraph: rewrite the "ability/knowledge ... type casting" text
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "This is an exercise that tests your ability/knowledge of collections, data associations, data modeling, serializations, type casting, and a number of other aspects of OOP/C#.NET.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "This is an exercise that tests your ability/knowledge with serialized data, collections, data associations, data modeling, type casting, and several other aspects of OOP/C#.NET.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Email paragraph: "please sent the follow up to" -> "please send the email to"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "If you are unable to reach them, please sent the follow up to ",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "If you are unable to reach them, please send the email to ",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. NOTE paragraph
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "NOTE: We will provide answers next to each question as a guide, but the answer should be programmatically determined in your work.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "NOTE: We will provide answers next to each most questions as a guide, but your answer should be programmatically determined in your work.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 4. Table cell text tweaks
# ---------------------------------------------------------------------------
# (search/replace text intentionally avoids the surrounding quote characters
# so Word's smart-quote autocorrect can't mangle the literal "B3" text)
$d.Content.Find.Execute(
    "breakdown of entities where the id contains",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "breakdown of entities where the EntityId contains",
    2) | Out-Null

$d.Content.Find.Execute(
    "Answer will not provided for this question.",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Answer will not be provided for this question.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5. Table vertical alignment: every data row (not the header row) switches
#    from "bottom" to "center".
# ---------------------------------------------------------------------------
$t = $d.Tables.Item(1)
for ($r = 2; $r -le $t.Rows.Count; $r++) {
    $row = $t.Rows.Item($r)
    for ($c = 1; $c -le $row.Cells.Count; $c++) {
        $row.Cells.Item($c).VerticalAlignment = 1
    }
}
